$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D column) values, preserving text formatting (e.g. trailing zeros)
# to match the source data which stores these as text, not numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.874.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.532.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.565"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.110"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.924.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.522.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.851.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0954"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0791"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0299"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.001.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.776.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.41"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) (E column) percentage values
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("E7").Value = "  -1.05%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("E15").Value = "  -3.39%  "
$ws.Range("E16").Value = "  -4.73%  "
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -4.90%  "
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("E33").Value = "  +7.55%  "
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("E37").Value = "  -4.46%  "
$ws.Range("E38").Value = "  -5.64%  "
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("E43").Value = "  +4.11%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("E47").Value = "  +3.27%  "
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("E51").Value = "  -1.42%  "
